$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, avoiding Excel's automatic
# number/date coercion for numeric-looking strings (e.g. "1.001", "27.105.68").
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "27.105.68"
Set-TextValue $ws.Range("E2") "  -2.73%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.716.44"
Set-TextValue $ws.Range("E3") "  -2.97%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  -0.15%  "

# Row 5
Set-TextValue $ws.Range("D5") "309.06"
Set-TextValue $ws.Range("E5") "  -5.77%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  -0.04%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4684"
Set-TextValue $ws.Range("E7") "  +4.45%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3431"
Set-TextValue $ws.Range("E8") "  -3.50%  "

# Row 9
Set-TextValue $ws.Range("D9") "42.08"
Set-TextValue $ws.Range("E9") "  +0.17%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.07265"
Set-TextValue $ws.Range("E10") "  -2.58%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.043"
Set-TextValue $ws.Range("E11") "  -5.01%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +0.05%  "

# Row 13
Set-TextValue $ws.Range("D13") "19.92"
Set-TextValue $ws.Range("E13") "  -5.00%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.856"
Set-TextValue $ws.Range("E14") "  -2.89%  "

# Row 15
Set-TextValue $ws.Range("D15") "1.714.56"
Set-TextValue $ws.Range("E15") "  -3.07%  "

# Row 16
Set-TextValue $ws.Range("D16") "6.880"
Set-TextValue $ws.Range("E16") "  -4.88%  "

# Row 17
Set-TextValue $ws.Range("D17") "88.66"
Set-TextValue $ws.Range("E17") "  -4.95%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.00001040"
Set-TextValue $ws.Range("E18") "  -1.73%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06363"
Set-TextValue $ws.Range("E19") "  -1.06%  "

# Row 20
Set-TextValue $ws.Range("D20") "1.000"
Set-TextValue $ws.Range("E20") "  -0.02%  "

# Row 21
Set-TextValue $ws.Range("D21") "16.54"
Set-TextValue $ws.Range("E21") "  -3.66%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.638"
Set-TextValue $ws.Range("E22") "  -2.40%  "

# Row 23
Set-TextValue $ws.Range("D23") "27.142.07"
Set-TextValue $ws.Range("E23") "  -2.79%  "

# Row 24
Set-TextValue $ws.Range("D24") "10.89"
Set-TextValue $ws.Range("E24") "  -3.76%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.138"
Set-TextValue $ws.Range("E25") "  +0.88%  "

# Row 26
Set-TextValue $ws.Range("D26") "156.79"
Set-TextValue $ws.Range("E26") "  -3.85%  "

# Row 27
Set-TextValue $ws.Range("D27") "19.43"
Set-TextValue $ws.Range("E27") "  -4.35%  "

# Row 28
Set-TextValue $ws.Range("D28") "1.908.00"
Set-TextValue $ws.Range("E28") "  -3.30%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.131"
Set-TextValue $ws.Range("E29") "  -1.62%  "

# Row 30
Set-TextValue $ws.Range("D30") "119.71"
Set-TextValue $ws.Range("E30") "  -4.31%  "

# Row 31
Set-TextValue $ws.Range("E31") "  -7.12%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.09150"
Set-TextValue $ws.Range("E32") "  -0.21%  "

# Row 33
Set-TextValue $ws.Range("D33") "3.588"
Set-TextValue $ws.Range("E33") "  -1.61%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.330"
Set-TextValue $ws.Range("E34") "  -4.87%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.02194"
Set-TextValue $ws.Range("E35") "  -4.33%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.05833"
Set-TextValue $ws.Range("E36") "  -4.80%  "

# Row 37
Set-TextValue $ws.Range("D37") "11.00"
Set-TextValue $ws.Range("E37") "  -7.44%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.1996"
Set-TextValue $ws.Range("E38") "  -5.02%  "

# Row 39
Set-TextValue $ws.Range("D39") "4.743"
Set-TextValue $ws.Range("E39") "  -4.49%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.393"
Set-TextValue $ws.Range("E40") "  -0.14%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.5910"
Set-TextValue $ws.Range("E41") "  -6.60%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -5.15%  "

# Row 43
Set-TextValue $ws.Range("D43") "7.467"
Set-TextValue $ws.Range("E43") "  -5.82%  "

# Row 44
Set-TextValue $ws.Range("D44") "12.64"
Set-TextValue $ws.Range("E44") "  -4.70%  "

# Row 45
Set-TextValue $ws.Range("B45") "Decentraland"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.5652"
Set-TextValue $ws.Range("E45") "  -3.87%  "

# Row 46
Set-TextValue $ws.Range("B46") "PancakeSwap"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D46") "3.562"
Set-TextValue $ws.Range("E46") "  -4.82%  "

# Row 47
Set-TextValue $ws.Range("D47") "118.55"
Set-TextValue $ws.Range("E47") "  -3.34%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.842"
Set-TextValue $ws.Range("E48") "  -5.83%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.06647"
Set-TextValue $ws.Range("E49") "  -3.71%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.086"
Set-TextValue $ws.Range("E50") "  -4.49%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.000"
Set-TextValue $ws.Range("E51") "  +0.06%  "
